$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4541.4165  # H64: 4453.615 -> 4541.4165
$ws.Cells.Item(64, 9).Value = 3624.5  # I64: 3579.6 -> 3624.5
$ws.Cells.Item(64, 11).Value = 3624.5  # K64: 3579.6 -> 3624.5
$ws.Cells.Item(64, 13).Value = -3376.5  # M64: -3331.6 -> -3376.5

$ws.Cells.Item(67, 8).Value = 4541.4165  # H67: 4453.615 -> 4541.4165
$ws.Cells.Item(67, 9).Value = 3624.5  # I67: 3579.6 -> 3624.5
$ws.Cells.Item(67, 11).Value = 3624.5  # K67: 3579.6 -> 3624.5
$ws.Cells.Item(67, 13).Value = -2766.5  # M67: -2721.6 -> -2766.5

$ws.Cells.Item(87, 8).Value = 129985  # H87: 129990 -> 129985
$ws.Cells.Item(87, 10).Value = 129985  # J87: 129990 -> 129985
$ws.Cells.Item(87, 12).Value = 129985  # L87: 129990 -> 129985
$ws.Cells.Item(87, 14).Value = -132481  # N87: -132486 -> -132481

$ws.Cells.Item(90, 8).Value = 129985  # H90: 129990 -> 129985
$ws.Cells.Item(90, 10).Value = 129985  # J90: 129990 -> 129985
$ws.Cells.Item(90, 12).Value = 389955  # L90: 389970 -> 389955
$ws.Cells.Item(90, 14).Value = -402435  # N90: -402450 -> -402435

$ws.Cells.Item(125, 8).Value = 2382  # H125: 2324.8333 -> 2382
$ws.Cells.Item(125, 9).Value = 1226.25  # I125: 1220.8 -> 1226.25
$ws.Cells.Item(125, 10).Value = 3042.4285  # J125: 3113.4285 -> 3042.4285
$ws.Cells.Item(125, 11).Value = 11036.25  # K125: 10987.2 -> 11036.25
$ws.Cells.Item(125, 12).Value = 27381.8565  # L125: 28020.8565 -> 27381.8565
$ws.Cells.Item(125, 13).Value = -8576.25  # M125: -8527.199999999999 -> -8576.25
$ws.Cells.Item(125, 14).Value = -32301.8565  # N125: -32940.8565 -> -32301.8565

$ws.Cells.Item(132, 8).Value = 700  # H132: 667.7368 -> 700
$ws.Cells.Item(132, 9).Value = 702.5192  # I132: 668.94543 -> 702.5192
$ws.Cells.Item(132, 11).Value = 2107.5576  # K132: 2006.83629 -> 2107.5576
$ws.Cells.Item(132, 13).Value = 422.4423999999999  # M132: 523.16371 -> 422.4423999999999

$ws.Cells.Item(137, 8).Value = 282509.94  # H137: 274885.34 -> 282509.94
$ws.Cells.Item(137, 9).Value = 387157.7  # I137: 372833.34 -> 387157.7
$ws.Cells.Item(137, 11).Value = 1161473.1  # K137: 1118500.02 -> 1161473.1
$ws.Cells.Item(137, 13).Value = -1158923.1  # M137: -1115950.02 -> -1158923.1

$ws.Cells.Item(141, 8).Value = 5300.75  # H141: 4524.5 -> 5300.75
$ws.Cells.Item(141, 9).Value = 4804.6665  # I141: 3450.6 -> 4804.6665
$ws.Cells.Item(141, 11).Value = 14413.9995  # K141: 10351.8 -> 14413.9995
$ws.Cells.Item(141, 13).Value = -9233.999500000002  # M141: -5171.799999999999 -> -9233.999500000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 8649.166999999999  # H63: 7700.143 -> 8649.166999999999
$ws.Cells.Item(63, 10).Value = 10000  # J63: 8667.666999999999 -> 10000
$ws.Cells.Item(63, 12).Value = 10000  # L63: 8667.666999999999 -> 10000
$ws.Cells.Item(63, 14).Value = -11372  # N63: -10039.667 -> -11372

$ws.Cells.Item(66, 8).Value = 8649.166999999999  # H66: 7700.143 -> 8649.166999999999
$ws.Cells.Item(66, 10).Value = 10000  # J66: 8667.666999999999 -> 10000
$ws.Cells.Item(66, 12).Value = 50000  # L66: 43338.335 -> 50000
$ws.Cells.Item(66, 14).Value = -56864  # N66: -50202.335 -> -56864

$ws.Cells.Item(74, 8).Value = 5325695  # H74: 5325567 -> 5325695
$ws.Cells.Item(74, 9).Value = 7814604.5  # I74: 7814416.5 -> 7814604.5
$ws.Cells.Item(74, 11).Value = 7814604.5  # K74: 7814416.5 -> 7814604.5
$ws.Cells.Item(74, 13).Value = -7813730.5  # M74: -7813542.5 -> -7813730.5

$ws.Cells.Item(77, 8).Value = 5325695  # H77: 5325567 -> 5325695
$ws.Cells.Item(77, 9).Value = 7814604.5  # I77: 7814416.5 -> 7814604.5
$ws.Cells.Item(77, 11).Value = 39073022.5  # K77: 39072082.5 -> 39073022.5
$ws.Cells.Item(77, 13).Value = -39068654.5  # M77: -39067714.5 -> -39068654.5

$ws.Cells.Item(110, 8).Value = 1848.75  # H110: 1884.4286 -> 1848.75
$ws.Cells.Item(110, 9).Value = 1848.75  # I110: 1881.8334 -> 1848.75
$ws.Cells.Item(110, 10).Value = 0  # J110: 1900 -> 0
$ws.Cells.Item(110, 11).Value = 1848.75  # K110: 1881.8334 -> 1848.75
$ws.Cells.Item(110, 12).Value = 0  # L110: 1900 -> 0
$ws.Cells.Item(110, 13).Value = 196.25  # M110: 163.1666 -> 196.25
$ws.Cells.Item(110, 14).ClearContents()  # N110: -5990 -> (removed)

$ws.Cells.Item(121, 8).Value = 82627.5  # H121: 81835 -> 82627.5
$ws.Cells.Item(121, 10).Value = 82627.5  # J121: 81835 -> 82627.5
$ws.Cells.Item(121, 12).Value = 82627.5  # L121: 81835 -> 82627.5
$ws.Cells.Item(121, 14).Value = -86121.5  # N121: -85329 -> -86121.5

$ws.Cells.Item(132, 8).Value = 7579806  # H132: 6806412.5 -> 7579806
$ws.Cells.Item(132, 9).Value = 11496773  # I132: 9806152 -> 11496773
$ws.Cells.Item(132, 11).Value = 34490319  # K132: 29418456 -> 34490319
$ws.Cells.Item(132, 13).Value = -34487789  # M132: -29415926 -> -34487789

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3674.0435  # H20: 2911.2258 -> 3674.0435
$ws.Cells.Item(20, 9).Value = 3365  # I20: 2754.8635 -> 3365
$ws.Cells.Item(20, 10).Value = 4786.6  # J20: 3293.4443 -> 4786.6
$ws.Cells.Item(20, 11).Value = 3365  # K20: 2754.8635 -> 3365
$ws.Cells.Item(20, 12).Value = 4786.6  # L20: 3293.4443 -> 4786.6
$ws.Cells.Item(20, 13).Value = -3118  # M20: -2507.8635 -> -3118
$ws.Cells.Item(20, 14).Value = -5280.6  # N20: -3787.4443 -> -5280.6

$ws.Cells.Item(23, 8).Value = 0  # H23: 2000 -> 0
$ws.Cells.Item(23, 10).Value = 0  # J23: 2000 -> 0
$ws.Cells.Item(23, 12).Value = 0  # L23: 2000 -> 0
$ws.Cells.Item(23, 14).ClearContents()  # N23: -2566 -> (removed)

$ws.Cells.Item(80, 8).Value = 1756  # H80: 1756.1111 -> 1756
$ws.Cells.Item(80, 9).Value = 249.5  # I80: 249.75 -> 249.5
$ws.Cells.Item(80, 11).Value = 249.5  # K80: 249.75 -> 249.5
$ws.Cells.Item(80, 13).Value = 748.5  # M80: 748.25 -> 748.5

$ws.Cells.Item(83, 8).Value = 1756  # H83: 1756.1111 -> 1756
$ws.Cells.Item(83, 9).Value = 249.5  # I83: 249.75 -> 249.5
$ws.Cells.Item(83, 11).Value = 1247.5  # K83: 1248.75 -> 1247.5
$ws.Cells.Item(83, 13).Value = 3744.5  # M83: 3743.25 -> 3744.5

$ws.Cells.Item(86, 8).Value = 2396.8262  # H86: 2390.88 -> 2396.8262
$ws.Cells.Item(86, 9).Value = 2311.7334  # I86: 2240.0557 -> 2311.7334
$ws.Cells.Item(86, 10).Value = 2556.375  # J86: 2778.7144 -> 2556.375
$ws.Cells.Item(86, 11).Value = 2311.7334  # K86: 2240.0557 -> 2311.7334
$ws.Cells.Item(86, 12).Value = 2556.375  # L86: 2778.7144 -> 2556.375
$ws.Cells.Item(86, 13).Value = -1188.7334  # M86: -1117.0557 -> -1188.7334
$ws.Cells.Item(86, 14).Value = -4802.375  # N86: -5024.7144 -> -4802.375

$ws.Cells.Item(89, 8).Value = 2396.8262  # H89: 2390.88 -> 2396.8262
$ws.Cells.Item(89, 9).Value = 2311.7334  # I89: 2240.0557 -> 2311.7334
$ws.Cells.Item(89, 10).Value = 2556.375  # J89: 2778.7144 -> 2556.375
$ws.Cells.Item(89, 11).Value = 11558.667  # K89: 11200.2785 -> 11558.667
$ws.Cells.Item(89, 12).Value = 12781.875  # L89: 13893.572 -> 12781.875
$ws.Cells.Item(89, 13).Value = -5942.667000000001  # M89: -5584.2785 -> -5942.667000000001
$ws.Cells.Item(89, 14).Value = -24013.875  # N89: -25125.572 -> -24013.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 745.125  # H16: 810.2857 -> 745.125
$ws.Cells.Item(16, 9).Value = 745.125  # I16: 810.2857 -> 745.125
$ws.Cells.Item(16, 11).Value = 745.125  # K16: 810.2857 -> 745.125
$ws.Cells.Item(16, 13).Value = -458.125  # M16: -523.2857 -> -458.125

$ws.Cells.Item(31, 8).Value = 1948078.4  # H31: 1798333.8 -> 1948078.4
$ws.Cells.Item(31, 9).Value = 2408.1667  # I31: 2264 -> 2408.1667
$ws.Cells.Item(31, 11).Value = 2408.1667  # K31: 2264 -> 2408.1667
$ws.Cells.Item(31, 13).Value = -2113.1667  # M31: -1969 -> -2113.1667

$ws.Cells.Item(34, 8).Value = 1948078.4  # H34: 1798333.8 -> 1948078.4
$ws.Cells.Item(34, 9).Value = 2408.1667  # I34: 2264 -> 2408.1667
$ws.Cells.Item(34, 11).Value = 2408.1667  # K34: 2264 -> 2408.1667
$ws.Cells.Item(34, 13).Value = -2206.1667  # M34: -2062 -> -2206.1667

$ws.Cells.Item(58, 8).Value = 6586.375  # H58: 7117.4287 -> 6586.375
$ws.Cells.Item(58, 10).Value = 7320.5  # J58: 8804.333000000001 -> 7320.5
$ws.Cells.Item(58, 12).Value = 7320.5  # L58: 8804.333000000001 -> 7320.5
$ws.Cells.Item(58, 14).Value = -7726.5  # N58: -9210.333000000001 -> -7726.5

$ws.Cells.Item(113, 8).Value = 745.125  # H113: 810.2857 -> 745.125
$ws.Cells.Item(113, 9).Value = 745.125  # I113: 810.2857 -> 745.125
$ws.Cells.Item(113, 11).Value = 745.125  # K113: 810.2857 -> 745.125
$ws.Cells.Item(113, 13).Value = 1424.875  # M113: 1359.7143 -> 1424.875

$ws.Cells.Item(132, 8).Value = 4491.0513  # H132: 4505.077 -> 4491.0513
$ws.Cells.Item(132, 9).Value = 1357.75  # I132: 1377.2858 -> 1357.75
$ws.Cells.Item(132, 11).Value = 4073.25  # K132: 4131.857400000001 -> 4073.25
$ws.Cells.Item(132, 13).Value = -1543.25  # M132: -1601.857400000001 -> -1543.25

$ws.Cells.Item(136, 8).Value = 6586.375  # H136: 7117.4287 -> 6586.375
$ws.Cells.Item(136, 10).Value = 7320.5  # J136: 8804.333000000001 -> 7320.5
$ws.Cells.Item(136, 12).Value = 21961.5  # L136: 26412.999 -> 21961.5
$ws.Cells.Item(136, 14).Value = -27061.5  # N136: -31512.999 -> -27061.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 165.10909  # H2: 163.41072 -> 165.10909
$ws.Cells.Item(2, 9).Value = 63.64  # I2: 63.884617 -> 63.64
$ws.Cells.Item(2, 11).Value = 381.84  # K2: 383.307702 -> 381.84
$ws.Cells.Item(2, 13).Value = -268.84  # M2: -270.307702 -> -268.84

$ws.Cells.Item(4, 8).Value = 10110436  # H4: 9515707 -> 10110436
$ws.Cells.Item(4, 9).Value = 8185904.5  # I4: 7503749.5 -> 8185904.5
$ws.Cells.Item(4, 11).Value = 24557713.5  # K4: 22511248.5 -> 24557713.5
$ws.Cells.Item(4, 13).Value = -24557601.5  # M4: -22511136.5 -> -24557601.5

$ws.Cells.Item(64, 8).Value = 4744  # H64: 4645 -> 4744
$ws.Cells.Item(64, 10).Value = 4766  # J64: 4649.125 -> 4766
$ws.Cells.Item(64, 12).Value = 14298  # L64: 13947.375 -> 14298
$ws.Cells.Item(64, 14).Value = -14838  # N64: -14487.375 -> -14838

$ws.Cells.Item(67, 8).Value = 4744  # H67: 4645 -> 4744
$ws.Cells.Item(67, 10).Value = 4766  # J67: 4649.125 -> 4766
$ws.Cells.Item(67, 12).Value = 14298  # L67: 13947.375 -> 14298
$ws.Cells.Item(67, 14).Value = -16170  # N67: -15819.375 -> -16170

$ws.Cells.Item(121, 8).Value = 1456.4783  # H121: 1454.2084 -> 1456.4783
$ws.Cells.Item(121, 9).Value = 512.4286  # I121: 584.1429000000001 -> 512.4286
$ws.Cells.Item(121, 10).Value = 1869.5  # J121: 1812.4706 -> 1869.5
$ws.Cells.Item(121, 11).Value = 1537.2858  # K121: 1752.4287 -> 1537.2858
$ws.Cells.Item(121, 12).Value = 5608.5  # L121: 5437.4118 -> 5608.5
$ws.Cells.Item(121, 13).Value = -227.2857999999999  # M121: -442.4287000000002 -> -227.2857999999999
$ws.Cells.Item(121, 14).Value = -8228.5  # N121: -8057.4118 -> -8228.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(44, 8).Value = 50000  # H44: 47500 -> 50000
$ws.Cells.Item(44, 9).Value = 50000  # I44: 47500 -> 50000
$ws.Cells.Item(44, 11).Value = 50000  # K44: 47500 -> 50000
$ws.Cells.Item(44, 13).Value = -49404  # M44: -46904 -> -49404

$ws.Cells.Item(47, 8).Value = 0  # H47: 32000 -> 0
$ws.Cells.Item(47, 10).Value = 0  # J47: 32000 -> 0
$ws.Cells.Item(47, 12).Value = 0  # L47: 32000 -> 0
$ws.Cells.Item(47, 14).ClearContents()  # N47: -33136 -> (removed)

$ws.Cells.Item(70, 8).Value = 6424.615  # H70: 6251.4287 -> 6424.615
$ws.Cells.Item(70, 9).Value = 5960.2856  # I70: 5715.25 -> 5960.2856
$ws.Cells.Item(70, 11).Value = 5960.2856  # K70: 5715.25 -> 5960.2856
$ws.Cells.Item(70, 13).Value = -5690.2856  # M70: -5445.25 -> -5690.2856

$ws.Cells.Item(73, 8).Value = 6424.615  # H73: 6251.4287 -> 6424.615
$ws.Cells.Item(73, 9).Value = 5960.2856  # I73: 5715.25 -> 5960.2856
$ws.Cells.Item(73, 11).Value = 5960.2856  # K73: 5715.25 -> 5960.2856
$ws.Cells.Item(73, 13).Value = -5024.2856  # M73: -4779.25 -> -5024.2856

$ws.Cells.Item(110, 8).Value = 102539  # H110: 102541.5 -> 102539
$ws.Cells.Item(110, 10).Value = 102539  # J110: 102541.5 -> 102539
$ws.Cells.Item(110, 12).Value = 102539  # L110: 102541.5 -> 102539
$ws.Cells.Item(110, 14).Value = -110719  # N110: -110721.5 -> -110719

$ws.Cells.Item(128, 8).Value = 117663.336  # H128: 117666.664 -> 117663.336
$ws.Cells.Item(128, 10).Value = 117663.336  # J128: 117666.664 -> 117663.336
$ws.Cells.Item(128, 12).Value = 117663.336  # L128: 117666.664 -> 117663.336
$ws.Cells.Item(128, 14).Value = -127623.336  # N128: -127626.664 -> -127623.336

$ws.Cells.Item(132, 8).Value = 52634976  # H132: 55559084 -> 52634976
$ws.Cells.Item(132, 9).Value = 71432340  # I132: 76927060 -> 71432340
$ws.Cells.Item(132, 11).Value = 214297020  # K132: 230781180 -> 214297020
$ws.Cells.Item(132, 13).Value = -214294490  # M132: -230778650 -> -214294490

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(75, 8).Value = 39993.332  # H75: 114000 -> 39993.332
$ws.Cells.Item(75, 9).Value = 19990  # I75: 0 -> 19990
$ws.Cells.Item(75, 10).Value = 49995  # J75: 114000 -> 49995
$ws.Cells.Item(75, 11).Value = 19990  # K75: 0 -> 19990
$ws.Cells.Item(75, 12).Value = 49995  # L75: 114000 -> 49995
$ws.Cells.Item(75, 14).Value = -51867  # N75: -115872 -> -51867
$ws.Cells.Item(75, 13).Value = -19054  # M75: None -> -19054

$ws.Cells.Item(78, 8).Value = 39993.332  # H78: 114000 -> 39993.332
$ws.Cells.Item(78, 9).Value = 19990  # I78: 0 -> 19990
$ws.Cells.Item(78, 10).Value = 49995  # J78: 114000 -> 49995
$ws.Cells.Item(78, 11).Value = 59970  # K78: 0 -> 59970
$ws.Cells.Item(78, 12).Value = 149985  # L78: 342000 -> 149985
$ws.Cells.Item(78, 14).Value = -159345  # N78: -351360 -> -159345
$ws.Cells.Item(78, 13).Value = -55290  # M78: None -> -55290

$ws.Cells.Item(100, 8).Value = 3418.7693  # H100: 3490.8333 -> 3418.7693
$ws.Cells.Item(100, 9).Value = 2677.7144  # I100: 2698.3333 -> 2677.7144
$ws.Cells.Item(100, 11).Value = 2677.7144  # K100: 2698.3333 -> 2677.7144
$ws.Cells.Item(100, 13).Value = -2136.7144  # M100: -2157.3333 -> -2136.7144

$ws.Cells.Item(127, 8).Value = 58242.715  # H127: 68331 -> 58242.715
$ws.Cells.Item(127, 10).Value = 58242.715  # J127: 68331 -> 58242.715
$ws.Cells.Item(127, 12).Value = 58242.715  # L127: 68331 -> 58242.715
$ws.Cells.Item(127, 14).Value = -68162.715  # N127: -78251 -> -68162.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 23134  # H52: 26012.334 -> 23134
$ws.Cells.Item(52, 9).Value = 17512.334  # I52: 19018.5 -> 17512.334
$ws.Cells.Item(52, 10).Value = 39999  # J52: 40000 -> 39999
$ws.Cells.Item(52, 11).Value = 17512.334  # K52: 19018.5 -> 17512.334
$ws.Cells.Item(52, 12).Value = 39999  # L52: 40000 -> 39999
$ws.Cells.Item(52, 13).Value = -17286.334  # M52: -18792.5 -> -17286.334
$ws.Cells.Item(52, 14).Value = -40451  # N52: -40452 -> -40451

$ws.Cells.Item(132, 8).Value = 3339853  # H132: 3092616 -> 3339853
$ws.Cells.Item(132, 9).Value = 6238.278  # I132: 5829.9 -> 6238.278
$ws.Cells.Item(132, 11).Value = 18714.834  # K132: 17489.7 -> 18714.834
$ws.Cells.Item(132, 13).Value = -16184.834  # M132: -14959.7 -> -16184.834

$ws.Cells.Item(133, 8).Value = 68333.336  # H133: 61428.75 -> 68333.336
$ws.Cells.Item(133, 10).Value = 68333.336  # J133: 61428.75 -> 68333.336
$ws.Cells.Item(133, 12).Value = 68333.336  # L133: 61428.75 -> 68333.336
$ws.Cells.Item(133, 14).Value = -78453.336  # N133: -71548.75 -> -78453.336
